$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3987.5862
$ws.Range("I40").Value = 3478.3125
$ws.Range("J40").Value = 4614.385
$ws.Range("K40").Value = 3478.3125
$ws.Range("L40").Value = 4614.385
$ws.Range("M40").Value = -3303.3125
$ws.Range("N40").Value = -4964.385
$ws.Range("H74").Value = 8350
$ws.Range("I74").Value = 7700
$ws.Range("K74").Value = 7700
$ws.Range("M74").Value = -6764
$ws.Range("H77").Value = 8350
$ws.Range("I77").Value = 7700
$ws.Range("K77").Value = 38500
$ws.Range("M77").Value = -33820
$ws.Range("H87").Value = 97499.5
$ws.Range("J87").Value = 97499.5
$ws.Range("L87").Value = 97499.5
$ws.Range("N87").Value = -99995.5
$ws.Range("H90").Value = 97499.5
$ws.Range("J90").Value = 97499.5
$ws.Range("L90").Value = 292498.5
$ws.Range("N90").Value = -304978.5
$ws.Range("H98").Value = 1088.25
$ws.Range("I98").Value = 858.8
$ws.Range("K98").Value = 858.8
$ws.Range("M98").Value = 639.2
$ws.Range("H122").Value = 1088.25
$ws.Range("I122").Value = 858.8
$ws.Range("K122").Value = 2576.4
$ws.Range("M122").Value = -126.3999999999996
$ws.Range("H125").Value = 2830
$ws.Range("I125").Value = 650
$ws.Range("J125").Value = 5010
$ws.Range("K125").Value = 5850
$ws.Range("L125").Value = 45090
$ws.Range("M125").Value = -3390
$ws.Range("N125").Value = -50010

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3564
$ws.Range("I2").Value = 3090.4614
$ws.Range("J2").Value = 5103
$ws.Range("K2").Value = 3090.4614
$ws.Range("L2").Value = 5103
$ws.Range("M2").Value = -2977.4614
$ws.Range("N2").Value = -5329
$ws.Range("H45").Value = 2609.7693
$ws.Range("I45").Value = 1387
$ws.Range("J45").Value = 5361
$ws.Range("K45").Value = 1387
$ws.Range("L45").Value = 5361
$ws.Range("M45").Value = -1010
$ws.Range("N45").Value = -6115
$ws.Range("H116").Value = 3564
$ws.Range("I116").Value = 3090.4614
$ws.Range("J116").Value = 5103
$ws.Range("K116").Value = 3090.4614
$ws.Range("L116").Value = 5103
$ws.Range("M116").Value = -796.4614000000001
$ws.Range("N116").Value = -9691
$ws.Range("H122").Value = 3315.25
$ws.Range("I122").Value = 1630.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 4891.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -2441.5
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 896.25
$ws.Range("I132").Value = 896.25
$ws.Range("K132").Value = 2688.75
$ws.Range("M132").Value = -158.75
$ws.Range("H135").Value = 100428.5
$ws.Range("J135").Value = 100428.5
$ws.Range("L135").Value = 100428.5
$ws.Range("N135").Value = -110568.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3564
$ws.Range("I3").Value = 3090.4614
$ws.Range("J3").Value = 5103
$ws.Range("K3").Value = 3090.4614
$ws.Range("L3").Value = 5103
$ws.Range("M3").Value = -2976.4614
$ws.Range("N3").Value = -5331

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2566.4443
$ws.Range("I22").Value = 1433
$ws.Range("J22").Value = 4833.3335
$ws.Range("K22").Value = 1433
$ws.Range("L22").Value = 4833.3335
$ws.Range("M22").Value = -1083
$ws.Range("N22").Value = -5533.3335
$ws.Range("H58").Value = 3302.4666
$ws.Range("I58").Value = 1645.7
$ws.Range("J58").Value = 6616
$ws.Range("K58").Value = 1645.7
$ws.Range("L58").Value = 6616
$ws.Range("M58").Value = -1442.7
$ws.Range("N58").Value = -7022
$ws.Range("H87").Value = 100000
$ws.Range("I87").Value = 100000
$ws.Range("K87").Value = 100000
$ws.Range("H90").Value = 100000
$ws.Range("I90").Value = 100000
$ws.Range("K90").Value = 300000
$ws.Range("H99").Value = 3619.9
$ws.Range("I99").Value = 2900
$ws.Range("J99").Value = 6499.5
$ws.Range("K99").Value = 2900
$ws.Range("L99").Value = 6499.5
$ws.Range("M99").Value = -1402
$ws.Range("N99").Value = -9495.5
$ws.Range("H126").Value = 3619.9
$ws.Range("I126").Value = 2900
$ws.Range("J126").Value = 6499.5
$ws.Range("K126").Value = 8700
$ws.Range("L126").Value = 19498.5
$ws.Range("M126").Value = -6230
$ws.Range("N126").Value = -24438.5
$ws.Range("H132").Value = 4687.375
$ws.Range("I132").Value = 4249.8335
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 12749.5005
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -10219.5005
$ws.Range("N132").Value = -23060
$ws.Range("H136").Value = 3302.4666
$ws.Range("I136").Value = 1645.7
$ws.Range("J136").Value = 6616
$ws.Range("K136").Value = 4937.1
$ws.Range("L136").Value = 19848
$ws.Range("M136").Value = -2387.1
$ws.Range("N136").Value = -24948
$ws.Range("M87").Value = -98814
$ws.Range("M90").Value = -294072

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 166717.33
$ws.Range("J4").Value = 86
$ws.Range("L4").Value = 258
$ws.Range("N4").Value = -482
$ws.Range("H80").Value = 4891
$ws.Range("J80").Value = 4825.6665
$ws.Range("L80").Value = 14476.9995
$ws.Range("N80").Value = -16348.9995
$ws.Range("H83").Value = 4891
$ws.Range("J83").Value = 4825.6665
$ws.Range("L83").Value = 43430.9985
$ws.Range("N83").Value = -52790.9985

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2931.9333
$ws.Range("I122").Value = 2315.3333
$ws.Range("J122").Value = 5398.3335
$ws.Range("K122").Value = 6945.999899999999
$ws.Range("L122").Value = 16195.0005
$ws.Range("M122").Value = -4495.999899999999
$ws.Range("N122").Value = -21095.0005
$ws.Range("H132").Value = 73545.36
$ws.Range("I132").Value = 85161.25
$ws.Range("J132").Value = 3850
$ws.Range("K132").Value = 255483.75
$ws.Range("L132").Value = 11550
$ws.Range("M132").Value = -252953.75
$ws.Range("N132").Value = -16610

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7568
$ws.Range("I7").Value = 6498.75
$ws.Range("J7").Value = 8993.666999999999
$ws.Range("K7").Value = 6498.75
$ws.Range("L7").Value = 8993.666999999999
$ws.Range("M7").Value = -6386.75
$ws.Range("N7").Value = -9217.666999999999
$ws.Range("H22").Value = 907.2222
$ws.Range("I22").Value = 909.2857
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 909.2857
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -614.2857
$ws.Range("N22").Value = -1490
$ws.Range("H27").Value = 907.2222
$ws.Range("I27").Value = 909.2857
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 909.2857
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -802.2857
$ws.Range("N27").Value = -1114
$ws.Range("H46").Value = 4834.1665
$ws.Range("J46").Value = 4999.5835
$ws.Range("L46").Value = 4999.5835
$ws.Range("N46").Value = -5375.5835
$ws.Range("H93").Value = 2460.5
$ws.Range("I93").Value = 2614
$ws.Range("K93").Value = 2614
$ws.Range("M93").Value = -1366
$ws.Range("H98").Value = 52225
$ws.Range("J98").Value = 52225
$ws.Range("L98").Value = 52225
$ws.Range("N98").Value = -58215
$ws.Range("H126").Value = 7568
$ws.Range("I126").Value = 6498.75
$ws.Range("J126").Value = 8993.666999999999
$ws.Range("K126").Value = 19496.25
$ws.Range("L126").Value = 26981.001
$ws.Range("M126").Value = -17026.25
$ws.Range("N126").Value = -31921.001
$ws.Range("H132").Value = 4498.6
$ws.Range("I132").Value = 4498.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13495.8
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10965.8
$ws.Range("N132").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1896.25
$ws.Range("I132").Value = 1952.8572
$ws.Range("K132").Value = 5858.571599999999
$ws.Range("M132").Value = -3328.571599999999
